$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = "What areas are inside 1000 foot of schools in El Cajon"
$ws.Range("D2").Value = "El Cajon"
$ws.Range("G2").Value = "amenity=school, amenity=kindergarten"
$ws.Range("J2").Value = "Buffer"
$ws.Range("K2").Value = "Overlay analysis"
$ws.Range("V2").Value = "data queries,buffer,overlay analysis"
$ws.Range("W2").Value = 10

# Row 3
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = "What areas are not park in Houston"
$ws.Range("D3").Value = "Houston"
$ws.Range("E3").Value = ""
$ws.Range("G3").Value = "leisure=park"
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("V3").Value = "data queries"
$ws.Range("W3").Value = 1

# Row 4
$ws.Range("A4").Value = 21
$ws.Range("B4").Value = 25
$ws.Range("C4").Value = "What areas are outside 250 meters of human settlement in the Cape Peninsula"
$ws.Range("D4").Value = "the Cape Peninsula"
$ws.Range("G4").Value = "residential=*"

# Row 5
$ws.Range("A5").Value = 23
$ws.Range("B5").Value = 27
$ws.Range("C5").Value = "What areas are outside 3000 meters of the rivers in Spain"
$ws.Range("D5").Value = "Spain"
$ws.Range("G5").Value = "waterway=river"

# Row 6
$ws.Range("A6").Value = 28
$ws.Range("B6").Value = 32
$ws.Range("C6").Value = "What areas are within 10 miles of current transmission lines with a voltage greater than 400 in Colorado"
$ws.Range("D6").Value = "Colorado"
$ws.Range("G6").Value = "power=line"

# Row 8
$ws.Range("A8").Value = 42
$ws.Range("B8").Value = 48
$ws.Range("C8").Value = "What areas are within 60 minutes of airports in Crook, Deschutes, and Jefferson county"
$ws.Range("D8").Value = "Crook, Deschutes, Jefferson county"
$ws.Range("G8").Value = " aeroway=*"
$ws.Range("J8").Value = "Network analysis"
$ws.Range("K8").Value = "classification"
$ws.Range("L8").Value = "Data queries"
$ws.Range("M8").Value = "Overlay analysis"
$ws.Range("V8").Value = "data queries,network analysis,classification,data queries,overlay analysis"
$ws.Range("W8").Value = 0

# Row 9
$ws.Range("A9").Value = 46
$ws.Range("B9").Value = 52
$ws.Range("C9").Value = "What areas are within a quarter mile of light rail stop in Gresham"
$ws.Range("D9").Value = "Gresham"
$ws.Range("G9").Value = "railway=tram_stop, light_rail"
$ws.Range("I9").Value = "Data queries"
$ws.Range("J9").Value = "Buffer"
$ws.Range("K9").Value = "Overlay analysis"
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("V9").Value = "data queries,buffer,overlay analysis"
$ws.Range("W9").Value = 10

# Row 10
$ws.Range("A10").Value = 54
$ws.Range("B10").Value = 61
$ws.Range("C10").Value = "What areas have an aspect larger than 45 degree and smaller than 135 degrees in the Cape Peninsula"
$ws.Range("D10").Value = "the Cape Peninsula"
$ws.Range("I10").Value = "Topography"
$ws.Range("J10").Value = "classification"
$ws.Range("K10").Value = "Data queries"
$ws.Range("L10").Value = "Data model conversion"
$ws.Range("M10").Value = "Overlay analysis"
$ws.Range("V10").Value = "topography,classification,data queries,data model conversion,overlay analysis"
$ws.Range("W10").Value = 17

# Row 11
$ws.Range("A11").Value = 60
$ws.Range("B11").Value = 68
$ws.Range("C11").Value = "What houses are for sale in urban areas in Utrecht"
$ws.Range("D11").Value = "Utrecht"
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = "boundary=urban, building=house"
$ws.Range("I11").Value = "Data queries"
$ws.Range("J11").Value = "Overlay analysis"
$ws.Range("K11").Value = "Data queries"
$ws.Range("V11").Value = "data queries,overlay analysis,data queries"
$ws.Range("W11").Value = 2

# Row 12
$ws.Range("A12").Value = 67
$ws.Range("B12").Value = 82
$ws.Range("C12").Value = "What is the density surface of temperature measurements in Oleander city"
$ws.Range("D12").Value = "Oleander city"
$ws.Range("G12").Value = " man_made=monitoring_station"

# Row 13
$ws.Range("A13").Value = 73
$ws.Range("B13").Value = 90
$ws.Range("C13").Value = "What is the Euclidean distance to the rivers in Crook, Deschutes, and Jefferson county"
$ws.Range("D13").Value = "Crook, Deschutes, Jefferson county"

# Row 14
$ws.Range("A14").Value = 80
$ws.Range("B14").Value = 97
$ws.Range("C14").Value = "What is the lung cancer mortality rate of white males for each city in the Western USA from 1970 to 1994"
$ws.Range("D14").Value = "the Western USA"
$ws.Range("F14").Value = " from 1970 to 1994"
$ws.Range("I14").Value = "data editing"
$ws.Range("J14").Value = "Data queries"
$ws.Range("V14").Value = "data editing,data queries"
$ws.Range("W14").Value = 9

# Row 15
$ws.Range("A15").Value = 82
$ws.Range("B15").Value = 99
$ws.Range("C15").Value = "What is the mean center of customers weighted by the number of transactions in Oleander city"
$ws.Range("D15").Value = "Oleander city"
$ws.Range("J15").Value = "Geostatistics  "
$ws.Range("K15").Value = ""
$ws.Range("V15").Value = "data queries,geostatistics  "
$ws.Range("W15").Value = 32

# Row 16
$ws.Range("A16").Value = 85
$ws.Range("B16").Value = 102
$ws.Range("C16").Value = "What is the mean center of the fire calls weighted by the priority in Fort Worth"
$ws.Range("D16").Value = "Fort Worth"
$ws.Range("G16").Value = ""
$ws.Range("J16").Value = "Overlay analysis"
$ws.Range("K16").Value = "Geostatistics  "
$ws.Range("L16").Value = ""
$ws.Range("V16").Value = "data queries,overlay analysis,geostatistics  "
$ws.Range("W16").Value = 35

# Row 17
$ws.Range("A17").Value = 87
$ws.Range("B17").Value = 105
$ws.Range("C17").Value = "What is the median people age for each census tract in Tarrant County"
$ws.Range("D17").Value = "Tarrant County"
$ws.Range("E17").Value = " Texas"
$ws.Range("G17").Value = ""
$ws.Range("I17").Value = "data editing"
$ws.Range("J17").Value = "Overlay analysis"
$ws.Range("K17").Value = "data editing"
$ws.Range("L17").Value = "Data queries"
$ws.Range("V17").Value = "data editing,overlay analysis,data editing,data queries"
$ws.Range("W17").Value = 26

# Row 18
$ws.Range("A18").Value = 93
$ws.Range("B18").Value = 111
$ws.Range("C18").Value = "Where are not conservation areas in UK"
$ws.Range("D18").Value = "UK"
$ws.Range("G18").Value = "landuse=conservation"

# Row 19
$ws.Range("A19").Value = 99
$ws.Range("B19").Value = 119
$ws.Range("C19").Value = "Where are the industrial areas in Utrecht"
$ws.Range("D19").Value = "Utrecht"
$ws.Range("G19").Value = "landuse=industrial"
$ws.Range("J19").Value = "Geometry measurement"
$ws.Range("K19").Value = "Data queries"
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("V19").Value = "data queries,geometry measurement,data queries"
$ws.Range("W19").Value = 8

# Row 20
$ws.Range("A20").Value = 110
$ws.Range("B20").Value = 132
$ws.Range("C20").Value = "Which shops are open at 6 pm in Happy Valley ski resort"
$ws.Range("D20").Value = "Happy Valley ski resort"
$ws.Range("G20").Value = "opening_hours=*"

# Row 21
$ws.Range("A21").Value = 111
$ws.Range("B21").Value = 133
$ws.Range("C21").Value = "Which vacant lots are within 1 mile of a freeway in Hillsboro"
$ws.Range("D21").Value = "Hillsboro"
$ws.Range("G21").Value = "abandoned:*=*, highway=motorway"
$ws.Range("J21").Value = "buffer"
$ws.Range("K21").Value = "Overlay analysis"
$ws.Range("L21").Value = "Data queries"
$ws.Range("V21").Value = "data queries,buffer,overlay analysis,data queries"
$ws.Range("W21").Value = 21

